$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data-row style (from B2) onto A2, which previously had
# the default style, to match the rest of row 2.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = "Rule import"
$ws.Range("B2").Value = "Custom Logistic"
$ws.Range("C2").Value = "Custom Logistic"
$ws.Range("D2").Value = "BA_CHD_Air"
$ws.Range("E2").Value = "SO_THC_OPS"
$ws.Range("F2").Value = 162789394
$ws.Range("F2").HorizontalAlignment = -4131
$ws.Range("G2").Value = "USAOMNIUSA"
$ws.Range("H2").Value = "Active"

$ws.Range("J10").Select()
